$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression values update
$ws.Range("B2").Value = 3662384476648496
$ws.Range("C2").Value = 3662384476648496
$ws.Range("D2").Value = 3662384476648496

# Row 3: RandomForestRegressor values update
$ws.Range("B3").Value = 3769118229239.652
$ws.Range("C3").Value = 62142447652384.75
$ws.Range("D3").Value = 684533435301580.6

# Row 4: rename GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03723385759680484
$ws.Range("C4").Value = 0.03614780946820455
$ws.Range("D4").Value = 331757495315105.8

# Row 5: rename AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 113826734634275.8
$ws.Range("C5").Value = 23753128016470.11
$ws.Range("D5").Value = 240557970665332.5
